# Update the header row (class labels) from "1-1", "1-2", "1-3" to "2-1", "2-2", "2-3"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "2-1"
$ws.Range("B1").Value = "2-2"
$ws.Range("C1").Value = "2-3"

# Move the active selection to C1, matching the saved view state in the target file
$ws.Range("C1").Select()
